# Natmi following Dr Hou advice
# Update the Pgf-Flt1 ligand-receptor pair table: recompute edge statistics
# across the full 3x3 cross-product of sending/target clusters (FAPs, sCs, ECs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds columns A-T for one data row (rows 2-10 of the sheet).
$rowData = @(
    @("FAPs", "Pgf", "Flt1", "FAPs", 2, 0.6666666666666666, 22.51188433333333, 67.535653, 0.7173237801266834, 0.7173237801266834, 3, 1, 61.629167, 184.887501, 0.8452417044501688, 0.8452417044501688, 1387.388679063684, 12486.49811157315, 0.606311974556916, 0.606311974556916),
    @("FAPs", "Pgf", "Flt1", "sCs", 2, 0.6666666666666666, 22.51188433333333, 67.535653, 0.7173237801266834, 0.7173237801266834, 3, 1, 1.357683666666667, 4.073051, 0.01862058035795759, 0.01862058035795759, 30.56401766525588, 275.076158987303, 0.01335698509052281, 0.01335698509052281),
    @("FAPs", "Pgf", "Flt1", "ECs", 2, 0.6666666666666666, 22.51188433333333, 67.535653, 0.7173237801266834, 0.7173237801266834, 3, 1, 9.926218666666665, 29.778656, 0.1361377151918736, 0.1361377151918736, 223.4578864913742, 2011.120978422368, 0.09765482047924456, 0.09765482047924458),
    @("sCs", "Pgf", "Flt1", "FAPs", 3, 1, 7.074492666666667, 21.223478, 0.2254232363222357, 0.2254232363222357, 3, 1, 61.629167, 184.887501, 0.8452417044501688, 0.8452417044501688, 435.9950899942754, 3923.955809948478, 0.1905371204916797, 0.1905371204916796),
    @("sCs", "Pgf", "Flt1", "sCs", 3, 1, 7.074492666666667, 21.223478, 0.2254232363222357, 0.2254232363222357, 3, 1, 1.357683666666667, 4.073051, 0.01862058035795759, 0.01862058035795759, 9.604923143486443, 86.444308291378, 0.004197511486489053, 0.004197511486489053),
    @("sCs", "Pgf", "Flt1", "ECs", 3, 1, 7.074492666666667, 21.223478, 0.2254232363222357, 0.2254232363222357, 3, 1, 9.926218666666665, 29.778656, 0.1361377151918736, 0.1361377151918736, 70.22296116506311, 632.006650485568, 0.03068860434406693, 0.03068860434406693),
    @("ECs", "Pgf", "Flt1", "FAPs", 3, 1, 1.796779333333333, 5.390338, 0.05725298355108089, 0.05725298355108089, 3, 1, 61.629167, 184.887501, 0.8452417044501688, 0.8452417044501688, 110.7340135961487, 996.606122365338, 0.04839260940157309, 0.04839260940157309),
    @("ECs", "Pgf", "Flt1", "sCs", 3, 1, 1.796779333333333, 5.390338, 0.05725298355108089, 0.05725298355108089, 3, 1, 1.357683666666667, 4.073051, 0.01862058035795759, 0.01862058035795759, 2.439457953470888, 21.955121581238, 0.001066083780945726, 0.001066083780945726),
    @("ECs", "Pgf", "Flt1", "ECs", 3, 1, 1.796779333333333, 5.390338, 0.05725298355108089, 0.05725298355108089, 3, 1, 9.926218666666665, 29.778656, 0.1361377151918736, 0.1361377151918736, 17.83522455841422, 160.517021025728, 0.007794290368562072, 0.007794290368562074),
)

$startRow = 2
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $rowData.Count; $i++) {
    $r = $startRow + $i
    $values = $rowData[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $ws.Range("$($columns[$j])$r").Value = $values[$j]
    }
}

